# Update crypto price / 1h-volume figures on the active worksheet to
# reflect the latest GitHub Actions scrape.
#
# Columns D (Price) and E (Volume(1h)) hold plain text in the original
# workbook (e.g. "519.42", "  -3.12%  ") even though many of the price
# strings parse as valid numbers. Assigning a numeric-looking string
# straight to .Value lets Excel auto-coerce the cell to a number (and
# subtly mangle the text, e.g. "520.39" -> 520.38999999999999), so the
# D column is temporarily forced to Text format while the new values are
# written, then restored to the default (Normal) style so no stray
# number formatting is left behind on the cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "58.149.16"
$ws.Range("E2").Value = "  -1.61%  "

$ws.Range("D3").Value = "2.473.19"
$ws.Range("E3").Value = "  -1.84%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "520.39"
$ws.Range("E5").Value = "  -2.85%  "

$ws.Range("D6").Value = "132.51"
$ws.Range("E6").Value = "  -3.72%  "

$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.10%  "

$ws.Range("E8").Value = "  -1.51%  "

$ws.Range("D9").Value = "0.0994"
$ws.Range("E9").Value = "  -1.75%  "

$ws.Range("E10").Value = "  -0.84%  "

$ws.Range("D11").Value = "5.38"
$ws.Range("E11").Value = "  +0.48%  "

$ws.Range("D12").Value = "0.343"
$ws.Range("E12").Value = "  -1.50%  "

$ws.Range("D13").Value = "2.910.18"
$ws.Range("E13").Value = "  -1.91%  "

$ws.Range("D14").Value = "58.072.13"
$ws.Range("E14").Value = "  -1.61%  "

$ws.Range("D15").Value = "22.09"
$ws.Range("E15").Value = "  -4.06%  "

$ws.Range("E16").Value = "  -1.97%  "

$ws.Range("D17").Value = "2.474.55"
$ws.Range("E17").Value = "  -1.79%  "

$ws.Range("D18").Value = "10.87"
$ws.Range("E18").Value = "  -2.11%  "

$ws.Range("E19").Value = "  -2.45%  "

$ws.Range("D20").Value = "320.78"
$ws.Range("E20").Value = "  -1.28%  "

$ws.Range("E21").Value = "  -0.14%  "

$ws.Range("D22").Value = "5.77"
$ws.Range("E22").Value = "  -3.03%  "

$ws.Range("D23").Value = "64.27"
$ws.Range("E23").Value = "  -2.33%  "

$ws.Range("D24").Value = "0.409"
$ws.Range("E24").Value = "  -3.19%  "

$ws.Range("E25").Value = "  -0.19%  "

$ws.Range("E26").Value = "  -3.73%  "

$ws.Range("D27").Value = "7.40"
$ws.Range("E27").Value = "  -3.13%  "

$ws.Range("D28").Value = "0.0₃0751"
$ws.Range("E28").Value = "  -2.38%  "

$ws.Range("D29").Value = "6.38"
$ws.Range("E29").Value = "  -4.67%  "

# Rows 30 / 31: Monero jumped above PancakeSwap in the ranking, so the
# two rows swap their Coin/Link/Price/Volume contents.
$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D30").Value = "166.94"
$ws.Range("E30").Value = "  +2.15%  "

$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "1.70"
$ws.Range("E31").Value = "  -4.65%  "

$ws.Range("E32").Value = "  -3.75%  "

$ws.Range("E33").Value = "  -0.04%  "

$ws.Range("E34").Value = "  -0.18%  "

$ws.Range("D35").Value = "18.15"
$ws.Range("E35").Value = "  -1.65%  "

$ws.Range("E36").Value = "  -10.00%  "

$ws.Range("D37").Value = "4.00"
$ws.Range("E37").Value = "  -3.00%  "

$ws.Range("E38").Value = "  -3.68%  "

$ws.Range("E39").Value = "  -2.86%  "

$ws.Range("D40").Value = "276.54"
$ws.Range("E40").Value = "  -3.37%  "

$ws.Range("D41").Value = "3.47"
$ws.Range("E41").Value = "  -4.43%  "

$ws.Range("D42").Value = "5.06"
$ws.Range("E42").Value = "  -2.68%  "

$ws.Range("D43").Value = "0.597"
$ws.Range("E43").Value = "  -1.54%  "

$ws.Range("D44").Value = "126.21"
$ws.Range("E44").Value = "  -4.60%  "

$ws.Range("E45").Value = "  -2.38%  "

$ws.Range("D46").Value = "0.0493"
$ws.Range("E46").Value = "  -3.33%  "

$ws.Range("D47").Value = "0.0215"
$ws.Range("E47").Value = "  -2.82%  "

$ws.Range("D48").Value = "17.15"
$ws.Range("E48").Value = "  -1.14%  "

$ws.Range("D49").Value = "1.737.64"
$ws.Range("E49").Value = "  -1.56%  "

$ws.Range("E50").Value = "  -1.81%  "

$ws.Range("D51").Value = "4.68"
$ws.Range("E51").Value = "  -1.72%  "

# Restore the default (unstyled) look for the price column now that the
# text values are safely in place.
$priceRange.Style = "Normal"
